# Update gh-pages to output generated at 2291077
# Adds a "Cover" column (J) to every sheet, and for the sheets that already
# contain event rows (展览 / 全部类型), fills in the cover-image URLs and
# refreshes the "想去人数" (F column) counters.

$wb = $excel.ActiveWorkbook

# Cover image URLs and updated "want to go" counts, keyed by row number.
$coverData = @{
    2 = @{ F = 15;   J = "//i0.hdslb.com/bfs/openplatform/202401/0PnysR0o1704703460388.jpeg" }
    3 = @{ F = 1784; J = "//i2.hdslb.com/bfs/openplatform/202312/l0hSA2KL1702521429527.jpeg" }
    4 = @{ F = 550;  J = "//i0.hdslb.com/bfs/openplatform/202312/QqiJ6HfK1702365336991.jpeg" }
    5 = @{ F = 1137; J = "//i1.hdslb.com/bfs/openplatform/202312/CIlf3jyZ1701747640038.jpeg" }
    6 = @{ F = 6011; J = "//i1.hdslb.com/bfs/openplatform/202311/YriBERx81701329557375.jpeg" }
    7 = @{ F = 141;  J = "//i1.hdslb.com/bfs/openplatform/202311/bv8DJewO1701071702232.jpeg" }
}

foreach ($ws in $wb.Worksheets) {
    # Every sheet gets a "Cover" header in column J.
    $ws.Range("J1").Value = "Cover"

    # Sheets that already have event rows (rows 2-7) get the cover URLs and
    # the refreshed "想去人数" counts.
    if ($ws.Cells.Item(7, 2).Text.Length -gt 0) {
        foreach ($row in $coverData.Keys) {
            $data = $coverData[$row]
            $ws.Cells.Item($row, 6).Value = $data.F
            $ws.Cells.Item($row, 10).Value = $data.J
        }
    }
}
